$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SS 2A")

$ws.Range("A3").Value = "2026-02-08 19:00:43"
$ws.Range("B3").Value = "Fatima Muhammed Gadaka"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "38"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 9
